$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.526.49'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.580.80'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D7").Value = '3.578.82'
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.22'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").Value = '4.193.55'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("E14").Value = '  -1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '3.581.35'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '65.571.89'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.586'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.44%  '
$ws.Range("D24").Value = '3.723.98'
$ws.Range("E24").Value = '  +0.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.17%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +29.55%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.98%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.997'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").Value = '3.581.32'
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.149'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '171.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0835'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.34%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("D50").Value = '2.470.60'
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0270'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.83%  '
